# HOC_Investor_Dashboard_Template.xlsx - Capital_Investment sheet
# Insert a new "Monthly Burn Rate (2027+)" row and update the burn-rate /
# runway figures, pushing the two trailing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: "Monthly Burn Rate" -> "Monthly Burn Rate (2026)" (value unchanged, note updated)
$ws.Range("A5").Value = "Monthly Burn Rate (2026)"
$ws.Range("B5").Value = 19239
$ws.Range("C5").Value = "2026 average (50% reduced rent): £5k rates + £12.8k rent + £1k service + £0.4k insurance"

# Row 6: new row - "Monthly Burn Rate (2027+)" (previously held "Runway (Months)")
$ws.Range("A6").Value = "Monthly Burn Rate (2027+)"
$ws.Range("B6").Value = 32078
$ws.Range("C6").Value = "2027+ average (full rent): £5k rates + £25.7k rent + £1k service + £0.4k insurance"

# Row 7: "Runway (Months)" shifts here (previously held "Next Major Expense"); value recalculated
$ws.Range("A7").Value = "Runway (Months)"
$ws.Range("B7").Value = 17
$ws.Range("C7").Value = "Until ~May 2027 (12 months 2026 + ~5 months 2027)"

# Row 8: "Next Major Expense" shifts here (previously held "Next Expense Description")
$ws.Range("A8").Value = "Next Major Expense"
$ws.Range("B8").Value = 244176.2
$ws.Range("C8").Value = "Due on landlord handover (Dec 19th)"

# Row 9: "Next Expense Description" - new trailing row
$ws.Range("A9").Value = "Next Expense Description"
$ws.Range("B9").Value = "December 2025 Total"
$ws.Range("C9").Value = "Rent deposit (7mo @ full rate), Q1 rent (50% rate), service charge, insurance, business rates, legal"

# Best-effort: mark the (now 9-row) table range as "number stored as text" so
# Excel doesn't flag the text-typed numeric-looking cells with a green
# triangle warning - mirrors the original A1:C8 ignoredError coverage.
$ws.Range("A1:C9").Errors.Item(9).Ignore = $true
